$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# row 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

# row 18
$ws.Range("H18").Value = 979.8333
$ws.Range("I18").Value = 1006.1818
$ws.Range("J18").Value = 690
$ws.Range("K18").Value = 1006.1818
$ws.Range("L18").Value = 690
$ws.Range("M18").Value = -722.1818
$ws.Range("N18").Value = -1258

# row 40
$ws.Range("H40").Value = 3000.5
$ws.Range("I40").Value = 2334
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2334
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2159
$ws.Range("N40").Value = -5350

# row 62
$ws.Range("H62").Value = 333335580
$ws.Range("I62").Value = 500001600
$ws.Range("J62").Value = 3499
$ws.Range("K62").Value = 500001600
$ws.Range("L62").Value = 3499
$ws.Range("M62").Value = -500000976
$ws.Range("N62").Value = -4747

# row 65
$ws.Range("H65").Value = 333335580
$ws.Range("I65").Value = 500001600
$ws.Range("J65").Value = 3499
$ws.Range("K65").Value = 2500008000
$ws.Range("L65").Value = 17495
$ws.Range("M65").Value = -2500004880
$ws.Range("N65").Value = -23735

# row 112
$ws.Range("H112").Value = 2590.8965
$ws.Range("I112").Value = 199
$ws.Range("J112").Value = 2676.3215
$ws.Range("K112").Value = 597
$ws.Range("L112").Value = 8028.9645
$ws.Range("M112").Value = 511
$ws.Range("N112").Value = -10244.9645

# row 113
$ws.Range("H113").Value = 2369.25
$ws.Range("I113").Value = 1188.5
$ws.Range("J113").Value = 3550
$ws.Range("K113").Value = 1188.5
$ws.Range("L113").Value = 3550
$ws.Range("M113").Value = 2065.5
$ws.Range("N113").Value = -10058

# row 135
$ws.Range("H135").Value = 6258927
$ws.Range("I135").Value = 11111849
$ws.Range("J135").Value = 19456.285
$ws.Range("K135").Value = 100006641
$ws.Range("L135").Value = 175106.565
$ws.Range("M135").Value = -100004106
$ws.Range("N135").Value = -180176.565

# row 141
$ws.Range("H141").Value = 1122.4286
$ws.Range("I141").Value = 971.4
$ws.Range("J141").Value = 1500
$ws.Range("K141").Value = 2914.2
$ws.Range("L141").Value = 4500
$ws.Range("M141").Value = 2265.8
$ws.Range("N141").Value = -14860

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# row 32
$ws.Range("H32").Value = 4084917
$ws.Range("I32").Value = 4654210
$ws.Range("J32").Value = 4983
$ws.Range("K32").Value = 4654210
$ws.Range("L32").Value = 4983
$ws.Range("M32").Value = -4653923
$ws.Range("N32").Value = -5557

# row 38
$ws.Range("H38").Value = 1900
$ws.Range("I38").Value = 1900
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1900
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -1433

# row 61
$ws.Range("H61").Value = 4538.625
$ws.Range("I61").Value = 3531.0417
$ws.Range("J61").Value = 7561.375
$ws.Range("K61").Value = 3531.0417
$ws.Range("L61").Value = 7561.375
$ws.Range("M61").Value = -3319.0417
$ws.Range("N61").Value = -7985.375

# row 74
$ws.Range("H74").Value = 4479.778
$ws.Range("I74").Value = 2239.6365
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 2239.6365
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -1365.6365
$ws.Range("N74").Value = -9748

# row 77
$ws.Range("H77").Value = 4479.778
$ws.Range("I77").Value = 2239.6365
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 11198.1825
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -6830.182500000001
$ws.Range("N77").Value = -48736

# row 132
$ws.Range("H132").Value = 5884.5938
$ws.Range("I132").Value = 3679.1155
$ws.Range("J132").Value = 15441.667
$ws.Range("K132").Value = 11037.3465
$ws.Range("L132").Value = 46325.001
$ws.Range("M132").Value = -8507.3465
$ws.Range("N132").Value = -51385.001

# row 136
$ws.Range("H136").Value = 4538.625
$ws.Range("I136").Value = 3531.0417
$ws.Range("J136").Value = 7561.375
$ws.Range("K136").Value = 10593.1251
$ws.Range("L136").Value = 22684.125
$ws.Range("M136").Value = -8043.125100000001
$ws.Range("N136").Value = -27784.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# row 86
$ws.Range("H86").Value = 2174.1177
$ws.Range("I86").Value = 1830.6666
$ws.Range("J86").Value = 4750
$ws.Range("K86").Value = 1830.6666
$ws.Range("L86").Value = 4750
$ws.Range("M86").Value = -707.6666
$ws.Range("N86").Value = -6996

# row 89
$ws.Range("H89").Value = 2174.1177
$ws.Range("I89").Value = 1830.6666
$ws.Range("J89").Value = 4750
$ws.Range("K89").Value = 9153.333000000001
$ws.Range("L89").Value = 23750
$ws.Range("M89").Value = -3537.333000000001
$ws.Range("N89").Value = -34982

# row 134
$ws.Range("H134").Value = 5167.659
$ws.Range("I134").Value = 2277.8
$ws.Range("J134").Value = 8970.105
$ws.Range("K134").Value = 6833.400000000001
$ws.Range("L134").Value = 26910.315
$ws.Range("M134").Value = -4298.400000000001
$ws.Range("N134").Value = -31980.315

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# row 22
$ws.Range("H22").Value = 443.58334
$ws.Range("I22").Value = 492.3
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 492.3
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -142.3
$ws.Range("N22").Value = -900

# row 31
$ws.Range("H31").Value = 5316.919
$ws.Range("I31").Value = 3398.5
$ws.Range("J31").Value = 6027.4443
$ws.Range("K31").Value = 3398.5
$ws.Range("L31").Value = 6027.4443
$ws.Range("M31").Value = -3103.5
$ws.Range("N31").Value = -6617.4443

# row 34
$ws.Range("H34").Value = 5316.919
$ws.Range("I34").Value = 3398.5
$ws.Range("J34").Value = 6027.4443
$ws.Range("K34").Value = 3398.5
$ws.Range("L34").Value = 6027.4443
$ws.Range("M34").Value = -3196.5
$ws.Range("N34").Value = -6431.4443

# row 41
$ws.Range("H41").Value = 14165.833
$ws.Range("I41").Value = 4999
$ws.Range("J41").Value = 60000
$ws.Range("K41").Value = 4999
$ws.Range("L41").Value = 60000
$ws.Range("M41").Value = -4571
$ws.Range("N41").Value = -60856

# row 86
$ws.Range("H86").Value = 8571.714
$ws.Range("I86").Value = 5998.8
$ws.Range("J86").Value = 15004
$ws.Range("K86").Value = 5998.8
$ws.Range("L86").Value = 15004
$ws.Range("M86").Value = -4875.8
$ws.Range("N86").Value = -17250

# row 89
$ws.Range("H89").Value = 8571.714
$ws.Range("I89").Value = 5998.8
$ws.Range("J89").Value = 15004
$ws.Range("K89").Value = 29994
$ws.Range("L89").Value = 75020
$ws.Range("M89").Value = -24378
$ws.Range("N89").Value = -86252

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# row 75
$ws.Range("H75").Value = 938.6667
$ws.Range("I75").Value = 938.6667
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2816.0001
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1818.0001
$ws.Range("N75").ClearContents()

# row 78
$ws.Range("H78").Value = 938.6667
$ws.Range("I78").Value = 938.6667
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 8448.0003
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -3456.0003
$ws.Range("N78").ClearContents()

# row 103
$ws.Range("H103").Value = 1930.6154
$ws.Range("I103").Value = 1909
$ws.Range("J103").Value = 2049.5
$ws.Range("K103").Value = 5727
$ws.Range("L103").Value = 6148.5
$ws.Range("M103").Value = -4848
$ws.Range("N103").Value = -7906.5

# row 105
$ws.Range("H105").Value = 12000
$ws.Range("I105").Value = 12000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 36000
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -33379

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()

# row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# row 55
$ws.Range("H55").Value = 574.2
$ws.Range("I55").Value = 292.27274
$ws.Range("J55").Value = 918.7778
$ws.Range("K55").Value = 292.27274
$ws.Range("L55").Value = 918.7778
$ws.Range("M55").Value = -119.27274
$ws.Range("N55").Value = -1264.7778

# row 58
$ws.Range("H58").Value = 7087
$ws.Range("I58").Value = 4830.6665
$ws.Range("J58").Value = 8440.799999999999
$ws.Range("K58").Value = 4830.6665
$ws.Range("L58").Value = 8440.799999999999
$ws.Range("M58").Value = -4570.6665
$ws.Range("N58").Value = -8960.799999999999

# row 122
$ws.Range("H122").Value = 20003912
$ws.Range("I122").Value = 29414594
$ws.Range("J122").Value = 6211.25
$ws.Range("K122").Value = 88243782
$ws.Range("L122").Value = 18633.75
$ws.Range("M122").Value = -88241332
$ws.Range("N122").Value = -23533.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# row 81
$ws.Range("H81").Value = 1701.7407
$ws.Range("I81").Value = 1382.5
$ws.Range("J81").Value = 10002
$ws.Range("K81").Value = 2765
$ws.Range("L81").Value = 20004
$ws.Range("M81").Value = -1704
$ws.Range("N81").Value = -22126

# row 84
$ws.Range("H84").Value = 1701.7407
$ws.Range("I84").Value = 1382.5
$ws.Range("J84").Value = 10002
$ws.Range("K84").Value = 13825
$ws.Range("L84").Value = 100020
$ws.Range("M84").Value = -8521
$ws.Range("N84").Value = -110628

# row 100
$ws.Range("H100").Value = 4976.091
$ws.Range("I100").Value = 6692.125
$ws.Range("J100").Value = 400
$ws.Range("K100").Value = 13384.25
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = -12843.25
$ws.Range("N100").Value = -1882

# row 126
$ws.Range("H126").Value = 5942941
$ws.Range("I126").Value = 7770000.5
$ws.Range("J126").Value = 4997.75
$ws.Range("K126").Value = 23310001.5
$ws.Range("L126").Value = 14993.25
$ws.Range("M126").Value = -23307531.5
$ws.Range("N126").Value = -19933.25

Write-Host "Applied Zalera_Profits updates"